# Add team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1, s="1":
# bold, centered, bordered) onto the three new header cells so they match
# the look of the rest of the header row.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-40) gets the same team record: 72 wins, 89 losses, 0 ties.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 72
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
